$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B category labels shift by two positions (a couple of new
# categories were inserted earlier in the master list used by the
# notebook, and the cached rows below were not regenerated, so every
# existing row's label slides down by two slots; "Thomas Hex" is also
# renamed to "Matthies Hex" along the way). ---
$ws.Range("B4").Value  = "Holden"
$ws.Range("B5").Value  = "Rizzie Spiral"
$ws.Range("B6").Value  = "RotRing OmegaMax-90"
$ws.Range("B7").Value  = "Equal Angle"
$ws.Range("B8").Value  = "Tilt Rotate"
$ws.Range("B9").Value  = "CLR"
$ws.Range("B10").Value = "Rizzie Hex"
$ws.Range("B11").Value = "Matthies Hex"
$ws.Range("B12").Value = "Tilt Rotate_Partial"
$ws.Range("B13").Value = "RotRing OmegaMax-60"
$ws.Range("B14").Value = "Equal Angle_Partial"
$ws.Range("B15").Value = "Rizzie Hex_Partial"
$ws.Range("B16").Value = "ND Single"
$ws.Range("B17").Value = "RD Single"
$ws.Range("B18").Value = "TD Single"
$ws.Range("B19").Value = "Morris Single"
$ws.Range("B20").Value = "Ring Perpendicular to ND"
$ws.Range("B21").Value = "Ring Perpendicular to RD"
$ws.Range("B22").Value = "Ring Perpendicular to TD"
$ws.Range("B23").Value = "OffsetFTD"
$ws.Range("B24").Value = "OffsetATD"
$ws.Range("B25").Value = "OffsetF45"
$ws.Range("B26").Value = "OffsetA45"
$ws.Range("B27").Value = "OffsetFRD"
$ws.Range("B28").Value = "OffsetARD"
$ws.Range("B29").Value = "Gaussian Quadrature"

# --- Add the two new data rows (30 and 31), mirroring the existing pattern ---
$newRows = @(
    @{ Row = 30; Id = 28; Label = "Michael-CCHex" },
    @{ Row = 31; Id = 29; Label = "Michael-SNHex" }
)

foreach ($item in $newRows) {
    $r = $item.Row

    # Copy formatting from the row above (keeps the same cell style index)
    $ws.Range("A" + ($r - 1)).Copy()
    $ws.Range("A" + $r).PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = $item.Id
    $ws.Cells.Item($r, 2).Value = $item.Label

    for ($c = 3; $c -le 23; $c++) {
        $ws.Cells.Item($r, $c).Value = 1
    }
}

$excel.CutCopyMode = 0
